# Applies the changes described in the commit diff:
#  - Switch the active/selected sheet from "Productos" to "Tipo_Afectacion"
#  - Update the selected cell on "Tipo_Afectacion" from C8 to C9
#  - Update values in Tipo_Afectacion!A2 (1 -> 5) and Tipo_Afectacion!A3 (2 -> 12)

$wb = $excel.ActiveWorkbook

# Update the data values on the "Tipo_Afectacion" sheet
$wsTipo = $wb.Worksheets.Item("Tipo_Afectacion")
$wsTipo.Range("A2").Value = 5
$wsTipo.Range("A3").Value = 12

# Make "Tipo_Afectacion" the active sheet (this updates workbook activeTab
# and moves tabSelected from "Productos" to "Tipo_Afectacion")
$wsTipo.Activate()

# Update the selection on the now-active sheet to C9
$wsTipo.Range("C9").Select()
